$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")
$ws.Range("A1").Value = $ws.Range("A1").Value
